# Scheduled runner update: refresh market-price-derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# "Omega_Profits" sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) with freshly
# pulled data. GSM has no changes this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1115.8572
$ws.Range("I33").Value = 482.2
$ws.Range("J33").Value = 2700
$ws.Range("K33").Value = 482.2
$ws.Range("L33").Value = 2700
$ws.Range("M33").Value = -253.2
$ws.Range("N33").Value = -3158
$ws.Range("H70").Value = 2668.3125
$ws.Range("J70").Value = 2726.6365
$ws.Range("L70").Value = 8179.9095
$ws.Range("N70").Value = -8719.9095
$ws.Range("H73").Value = 2668.3125
$ws.Range("J73").Value = 2726.6365
$ws.Range("L73").Value = 8179.9095
$ws.Range("N73").Value = -10051.9095
$ws.Range("H88").Value = 2944.55
$ws.Range("J88").Value = 3211.8125
$ws.Range("L88").Value = 3211.8125
$ws.Range("N88").Value = -4023.8125
$ws.Range("H91").Value = 2944.55
$ws.Range("J91").Value = 3211.8125
$ws.Range("L91").Value = 3211.8125
$ws.Range("N91").Value = -6019.8125
$ws.Range("H101").Value = 1782.5834
$ws.Range("I101").Value = 1490.091
$ws.Range("K101").Value = 4470.272999999999
$ws.Range("M101").Value = -2848.272999999999
$ws.Range("H112").Value = 6266.65
$ws.Range("J112").Value = 6266.65
$ws.Range("L112").Value = 18799.95
$ws.Range("N112").Value = -21015.95
$ws.Range("H132").Value = 2822.0862
$ws.Range("I132").Value = 2596.1455
$ws.Range("K132").Value = 7788.4365
$ws.Range("M132").Value = -5258.4365
$ws.Range("H138").Value = 3327.2307
$ws.Range("I138").Value = 2983.15
$ws.Range("J138").Value = 3424.155
$ws.Range("K138").Value = 8949.450000000001
$ws.Range("L138").Value = 10272.465
$ws.Range("M138").Value = -3809.450000000001
$ws.Range("N138").Value = -20552.465
$ws.Range("H139").Value = 89999.8
$ws.Range("J139").Value = 89999.8
$ws.Range("L139").Value = 89999.8
$ws.Range("N139").Value = -100279.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13150.737
$ws.Range("I32").Value = 8332.058000000001
$ws.Range("K32").Value = 8332.058000000001
$ws.Range("M32").Value = -8045.058000000001
$ws.Range("H88").Value = 3028.7273
$ws.Range("J88").Value = 2966.25
$ws.Range("L88").Value = 2966.25
$ws.Range("N88").Value = -3778.25
$ws.Range("H91").Value = 3028.7273
$ws.Range("J91").Value = 2966.25
$ws.Range("L91").Value = 2966.25
$ws.Range("N91").Value = -5774.25
$ws.Range("H97").Value = 460.30768
$ws.Range("I97").Value = 460.30768
$ws.Range("K97").Value = 460.30768
$ws.Range("M97").Value = 35.69232
$ws.Range("H98").Value = 85999.82000000001
$ws.Range("J98").Value = 85999.82000000001
$ws.Range("L98").Value = 85999.82000000001
$ws.Range("N98").Value = -91989.82000000001
$ws.Range("H132").Value = 2818.0952
$ws.Range("I132").Value = 2093.257
$ws.Range("K132").Value = 6279.771000000001
$ws.Range("M132").Value = -3749.771000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 89999.836
$ws.Range("J57").Value = 89999.836
$ws.Range("L57").Value = 89999.836
$ws.Range("N57").Value = -91439.836
$ws.Range("H60").Value = 35994.6
$ws.Range("J60").Value = 35994.6
$ws.Range("L60").Value = 35994.6
$ws.Range("N60").Value = -37192.6
$ws.Range("H94").Value = 1050.0667
$ws.Range("I94").Value = 1028
$ws.Range("K94").Value = 1028
$ws.Range("M94").Value = -577
$ws.Range("H100").Value = 36814
$ws.Range("J100").Value = 36814
$ws.Range("L100").Value = 36814
$ws.Range("N100").Value = -38978
$ws.Range("H133").Value = 89990.91
$ws.Range("J133").Value = 89990.91
$ws.Range("L133").Value = 89990.91
$ws.Range("N133").Value = -100110.91
$ws.Range("H136").Value = 89999.836
$ws.Range("J136").Value = 89999.836
$ws.Range("L136").Value = 89999.836
$ws.Range("N136").Value = -100199.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6200.5293
$ws.Range("I31").Value = 6968.273
$ws.Range("J31").Value = 5618.1035
$ws.Range("K31").Value = 6968.273
$ws.Range("L31").Value = 5618.1035
$ws.Range("M31").Value = -6673.273
$ws.Range("N31").Value = -6208.1035
$ws.Range("H34").Value = 6200.5293
$ws.Range("I34").Value = 6968.273
$ws.Range("J34").Value = 5618.1035
$ws.Range("K34").Value = 6968.273
$ws.Range("L34").Value = 5618.1035
$ws.Range("M34").Value = -6766.273
$ws.Range("N34").Value = -6022.1035
$ws.Range("H50").Value = 34999
$ws.Range("J50").Value = 34999
$ws.Range("L50").Value = 34999
$ws.Range("N50").Value = -36249
$ws.Range("H52").Value = 31224.75
$ws.Range("I52").Value = 27500
$ws.Range("J52").Value = 34949.5
$ws.Range("K52").Value = 27500
$ws.Range("L52").Value = 34949.5
$ws.Range("M52").Value = -27206
$ws.Range("N52").Value = -35537.5
$ws.Range("H60").Value = 29757
$ws.Range("J60").Value = 29757
$ws.Range("L60").Value = 29757
$ws.Range("N60").Value = -30779
$ws.Range("H132").Value = 1005.3182
$ws.Range("I132").Value = 889.5789
$ws.Range("K132").Value = 2668.7367
$ws.Range("M132").Value = -138.7366999999999
$ws.Range("H133").Value = 89416.164
$ws.Range("J133").Value = 89416.164
$ws.Range("L133").Value = 89416.164
$ws.Range("N133").Value = -94476.164
$ws.Range("H134").Value = 1716.5
$ws.Range("I134").Value = 1605.16
$ws.Range("K134").Value = 4815.48
$ws.Range("M134").Value = -2280.48
$ws.Range("H137").Value = 86909.69
$ws.Range("I137").Value = 76498.5
$ws.Range("J137").Value = 88397
$ws.Range("K137").Value = 76498.5
$ws.Range("L137").Value = 88397
$ws.Range("M137").Value = -71398.5
$ws.Range("N137").Value = -98597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2404.3125
$ws.Range("I129").Value = 1598
$ws.Range("J129").Value = 2770.818
$ws.Range("K129").Value = 4794
$ws.Range("L129").Value = 8312.454000000002
$ws.Range("M129").Value = 206
$ws.Range("N129").Value = -18312.454
$ws.Range("H131").Value = 3463.4546
$ws.Range("I131").Value = 998.5
$ws.Range("J131").Value = 4872
$ws.Range("K131").Value = 2995.5
$ws.Range("L131").Value = 14616
$ws.Range("M131").Value = 2044.5
$ws.Range("N131").Value = -24696
$ws.Range("H141").Value = 8251.308000000001
$ws.Range("I141").Value = 6207.778
$ws.Range("K141").Value = 18623.334
$ws.Range("M141").Value = -13443.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2745.923
$ws.Range("I61").Value = 2641.4167
$ws.Range("K61").Value = 2641.4167
$ws.Range("M61").Value = -2439.4167
$ws.Range("H68").Value = 4179.4
$ws.Range("J68").Value = 3633
$ws.Range("L68").Value = 3633
$ws.Range("N68").Value = -5131
$ws.Range("H71").Value = 4179.4
$ws.Range("J71").Value = 3633
$ws.Range("L71").Value = 18165
$ws.Range("N71").Value = -25653
$ws.Range("H113").Value = 2745.923
$ws.Range("I113").Value = 2641.4167
$ws.Range("K113").Value = 2641.4167
$ws.Range("M113").Value = -471.4167000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101248
$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306240
$ws.Range("H132").Value = 4502
$ws.Range("I132").Value = 4005.3667
$ws.Range("K132").Value = 12016.1001
$ws.Range("M132").Value = -9486.1001
